# Updated the users list and fixtures list with new fixtures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop the leftover "UNI_*" hidden defined names (exported artifact)
# ---------------------------------------------------------------------
while ($wb.Names.Count() -gt 0) {
    $wb.Names.Item(1).Delete()
}

# ---------------------------------------------------------------------
# 2) Fill in the Round-of-16 results that had just been played
#    (columns J = Home_Score, K = Away_Score)
# ---------------------------------------------------------------------
$ws.Range("J38").Value = 2
$ws.Range("K38").Value = 0

$ws.Range("J39").Value = 2
$ws.Range("K39").Value = 0

$ws.Range("J40").Value = 1
$ws.Range("K40").Value = 1

$ws.Range("J41").Value = 4
$ws.Range("K41").Value = 1

$ws.Range("J42").Value = 1
$ws.Range("K42").Value = 0

$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0

$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 3

$ws.Range("J45").Value = 1
$ws.Range("K45").Value = 2

# ---------------------------------------------------------------------
# 3) Append the newly scheduled quarter-final fixtures (rows 46-49)
# ---------------------------------------------------------------------
$ws.Range("A46").Value = "Fri"
$ws.Range("B46").Value = "Jul 05, 2024"
$ws.Range("C46").Value = "18:00:00"
$ws.Range("D46").Value = "Spain"
$ws.Range("G46").Value = "Germany"
$ws.Range("H46").Value = "Stuttgart"

$ws.Range("A47").Value = "Fri"
$ws.Range("B47").Value = "Jul 05, 2024"
$ws.Range("C47").Value = "21:00:00"
$ws.Range("D47").Value = "Portugal"
$ws.Range("G47").Value = "France"
$ws.Range("H47").Value = "Hamburg"

$ws.Range("A48").Value = "Sat"
$ws.Range("B48").Value = "Jul 06, 2024"
$ws.Range("C48").Value = "18:00:00"
$ws.Range("D48").Value = "England"
$ws.Range("G48").Value = "Switzerland"
$ws.Range("H48").Value = "Düsseldorf"

$ws.Range("A49").Value = "Sat"
$ws.Range("B49").Value = "Jul 06, 2024"
$ws.Range("C49").Value = "21:00:00"
$ws.Range("D49").Value = "Netherlands"
$ws.Range("G49").Value = "Turkey"
$ws.Range("H49").Value = "Berlin"

# ---------------------------------------------------------------------
# 4) Restore the viewport/selection to what the author left it at
# ---------------------------------------------------------------------
$excel.Goto($ws.Range("A16"))
$ws.Range("J48").Select()
